$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$t.Cell(1,1).Range.Text = "29 x 18" + $nl + "  1    8" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"
$t.Cell(1,2).Range.Text = "15 x 81" + $nl + "  8    1" + $nl + "  ----" + $nl + "1|    |" + $nl + "5|    |"
$t.Cell(1,3).Range.Text = "69 x 72" + $nl + "  7    2" + $nl + "  ----" + $nl + "6|    |" + $nl + "9|    |"
$t.Cell(2,1).Range.Text = "87 x 10" + $nl + "  1    0" + $nl + "  ----" + $nl + "8|    |" + $nl + "7|    |"
$t.Cell(2,2).Range.Text = "42 x 71" + $nl + "  7    1" + $nl + "  ----" + $nl + "4|    |" + $nl + "2|    |"
$t.Cell(2,3).Range.Text = "73 x 18" + $nl + "  1    8" + $nl + "  ----" + $nl + "7|    |" + $nl + "3|    |"
$t.Cell(3,1).Range.Text = "39 x 82" + $nl + "  8    2" + $nl + "  ----" + $nl + "3|    |" + $nl + "9|    |"
$t.Cell(3,2).Range.Text = "25 x 50" + $nl + "  5    0" + $nl + "  ----" + $nl + "2|    |" + $nl + "5|    |"
$t.Cell(3,3).Range.Text = "76 x 69" + $nl + "  6    9" + $nl + "  ----" + $nl + "7|    |" + $nl + "6|    |"
$t.Cell(4,1).Range.Text = "26 x 80" + $nl + "  8    0" + $nl + "  ----" + $nl + "2|    |" + $nl + "6|    |"
$t.Cell(4,2).Range.Text = "80 x 93" + $nl + "  9    3" + $nl + "  ----" + $nl + "8|    |" + $nl + "0|    |"
$t.Cell(4,3).Range.Text = "83 x 43" + $nl + "  4    3" + $nl + "  ----" + $nl + "8|    |" + $nl + "3|    |"
$t.Cell(5,1).Range.Text = "47 x 83" + $nl + "  8    3" + $nl + "  ----" + $nl + "4|    |" + $nl + "7|    |"
$t.Cell(5,2).Range.Text = "42 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "4|    |" + $nl + "2|    |"
$t.Cell(5,3).Range.Text = "93 x 35" + $nl + "  3    5" + $nl + "  ----" + $nl + "9|    |" + $nl + "0|    |"
